# Sprint 4 Backlog - Burndown: mark "Remove ingredients used to cook recipe"
# tasks as completed for the web functionality (row 13) and assign the
# desktop companion tasks (rows 17 & 19) to Matthew.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 - "Complete functionality to remove ingredients used to cook
# recipe from pantry (web)" : record 1 hour actual time, completed by
# Destiny, and 0 hours remaining.
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "Destiny"
$ws.Range("I13").Value = 0

# Row 17 - "Create UI to view shared recipes (desktop)" : 0.5 hour
# estimate assigned to Matthew.
$ws.Range("C17").Value = 0.5
$ws.Range("D17").Value = "Matthew"

# Row 19 - "Create UI to view shared recipes (web)" : 0.5 hour estimate
# assigned to Matthew.
$ws.Range("C19").Value = 0.5
$ws.Range("D19").Value = "Matthew"

# Update the active selection / scroll position left by the author after
# making these edits.
$ws.Range("E12").Select()
